# Fruta / hortaliza, semanal
# Insert a new weekly data row at the top of the data block (row 73) and
# shift the existing rows 73:90 down to 74:91 (the oldest row falls through
# to the new last row, 91).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above row 73; this pushes rows 73-90 down to 74-91
# and carries the row's number formatting along automatically.
$ws.Rows("73:73").Insert()

# Populate the newly inserted row 73 with this week's record.
$ws.Range("A73").Value = 7
$ws.Range("B73").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C73").Value = "Ñuble"
$ws.Range("D73").Value2 = 45218
$ws.Range("E73").Value = 16
$ws.Range("F73").Value = 100112026
$ws.Range("G73").Value = "Haba"
$ws.Range("H73").Value = "Sin especificar"
$ws.Range("I73").Value = "Primera"
$ws.Range("J73").Value = 80
$ws.Range("K73").Value = 15000
$ws.Range("L73").Value = 15000
$ws.Range("M73").Value = 15000
$ws.Range("N73").Value = "$/saco 25 kilos"
$ws.Range("O73").Value = "Provincia de Diguillín"
$ws.Range("P73").Value = 600
$ws.Range("Q73").Value = 25
$ws.Range("R73").Value = "Hortaliza"
